# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data refresh to Sheets/Faerie_Profits.xlsx
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N) for affected rows

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 435.5
$ws.Range("I55").Value = 48.5
$ws.Range("K55").Value = 48.5
$ws.Range("M55").Value = 165.5
$ws.Range("H62").Value = 6530.222
$ws.Range("I62").Value = 6236.4165
$ws.Range("K62").Value = 6236.4165
$ws.Range("M62").Value = -5612.4165
$ws.Range("H64").Value = 8388.429
$ws.Range("I64").Value = 5749.5
$ws.Range("J64").Value = 9444
$ws.Range("K64").Value = 5749.5
$ws.Range("L64").Value = 9444
$ws.Range("M64").Value = -5501.5
$ws.Range("N64").Value = -9940
$ws.Range("H65").Value = 6530.222
$ws.Range("I65").Value = 6236.4165
$ws.Range("K65").Value = 31182.0825
$ws.Range("M65").Value = -28062.0825
$ws.Range("H67").Value = 8388.429
$ws.Range("I67").Value = 5749.5
$ws.Range("J67").Value = 9444
$ws.Range("K67").Value = 5749.5
$ws.Range("L67").Value = 9444
$ws.Range("M67").Value = -4891.5
$ws.Range("N67").Value = -11160
$ws.Range("H92").Value = 825.4815
$ws.Range("I92").Value = 709.4231
$ws.Range("K92").Value = 709.4231
$ws.Range("M92").Value = 538.5769
$ws.Range("H96").Value = 1215.6
$ws.Range("I96").Value = 1270.6666
$ws.Range("J96").Value = 720
$ws.Range("K96").Value = 3811.9998
$ws.Range("L96").Value = 2160
$ws.Range("M96").Value = -2438.9998
$ws.Range("N96").Value = -4906
$ws.Range("H97").Value = 2598.4
$ws.Range("I97").Value = 2249.5
$ws.Range("J97").Value = 2831
$ws.Range("K97").Value = 6748.5
$ws.Range("L97").Value = 8493
$ws.Range("M97").Value = -6252.5
$ws.Range("N97").Value = -9485
$ws.Range("H98").Value = 5498289.5
$ws.Range("I98").Value = 7144293.5
$ws.Range("K98").Value = 7144293.5
$ws.Range("M98").Value = -7142795.5
$ws.Range("H99").Value = 19608278
$ws.Range("I99").Value = 25641344
$ws.Range("K99").Value = 76924032
$ws.Range("M99").Value = -76922534
$ws.Range("H100").Value = 5915.4546
$ws.Range("I100").Value = 1064.2
$ws.Range("K100").Value = 1064.2
$ws.Range("M100").Value = -523.2
$ws.Range("H101").Value = 142857570
$ws.Range("I101").Value = 20408652
$ws.Range("K101").Value = 61225956
$ws.Range("M101").Value = -61224334
$ws.Range("H103").Value = 366.92856
$ws.Range("I103").Value = 230.25
$ws.Range("K103").Value = 690.75
$ws.Range("M103").Value = -104.75
$ws.Range("H116").Value = 3549.6924
$ws.Range("I116").Value = 3595.5
$ws.Range("K116").Value = 3595.5
$ws.Range("M116").Value = -153.5
$ws.Range("H118").Value = 2434.9
$ws.Range("I118").Value = 2434.9
$ws.Range("K118").Value = 7304.700000000001
$ws.Range("M118").Value = -5647.700000000001
$ws.Range("H122").Value = 5498289.5
$ws.Range("I122").Value = 7144293.5
$ws.Range("K122").Value = 21432880.5
$ws.Range("M122").Value = -21430430.5
$ws.Range("H138").Value = 334551.53
$ws.Range("I138").Value = 38000.297
$ws.Range("J138").Value = 1669032.1
$ws.Range("K138").Value = 114000.891
$ws.Range("L138").Value = 5007096.300000001
$ws.Range("M138").Value = -108860.891
$ws.Range("N138").Value = -5017376.300000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2823.9
$ws.Range("I45").Value = 2405.5
$ws.Range("K45").Value = 2405.5
$ws.Range("M45").Value = -2028.5
$ws.Range("H102").Value = 2470.9565
$ws.Range("I102").Value = 1509.5
$ws.Range("K102").Value = 1509.5
$ws.Range("M102").Value = 112.5
$ws.Range("H132").Value = 1922.1351
$ws.Range("I132").Value = 1508
$ws.Range("K132").Value = 4524
$ws.Range("M132").Value = -1994

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 975.875
$ws.Range("I94").Value = 217.83333
$ws.Range("K94").Value = 217.83333
$ws.Range("M94").Value = 233.16667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1999.5555
$ws.Range("I16").Value = 1574.5
$ws.Range("K16").Value = 1574.5
$ws.Range("M16").Value = -1287.5
$ws.Range("H58").Value = 3268.9
$ws.Range("I58").Value = 3737.8
$ws.Range("J58").Value = 2800
$ws.Range("K58").Value = 3737.8
$ws.Range("L58").Value = 2800
$ws.Range("M58").Value = -3534.8
$ws.Range("N58").Value = -3206
$ws.Range("H94").Value = 2563
$ws.Range("J94").Value = 2514.5557
$ws.Range("L94").Value = 2514.5557
$ws.Range("N94").Value = -3416.5557
$ws.Range("H113").Value = 1999.5555
$ws.Range("I113").Value = 1574.5
$ws.Range("K113").Value = 1574.5
$ws.Range("M113").Value = 595.5
$ws.Range("H136").Value = 3268.9
$ws.Range("I136").Value = 3737.8
$ws.Range("J136").Value = 2800
$ws.Range("K136").Value = 11213.4
$ws.Range("L136").Value = 8400
$ws.Range("M136").Value = -8663.400000000001
$ws.Range("N136").Value = -13500
$ws.Range("H138").Value = 82945
$ws.Range("J138").Value = 82945
$ws.Range("L138").Value = 82945
$ws.Range("N138").Value = -93225
$ws.Range("H139").Value = 100740
$ws.Range("J139").Value = 100740
$ws.Range("L139").Value = 100740
$ws.Range("N139").Value = -111020
$ws.Range("H140").Value = 112829.945
$ws.Range("I140").Value = 117857.07
$ws.Range("K140").Value = 117857.07
$ws.Range("M140").Value = -112677.07
$ws.Range("H141").Value = 159221.17
$ws.Range("J141").Value = 159221.17
$ws.Range("L141").Value = 159221.17
$ws.Range("N141").Value = -169581.17

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 5499.8335
$ws.Range("J104").Value = 5499.8335
$ws.Range("L104").Value = 16499.5005
$ws.Range("N104").Value = -21741.5005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2973.3845
$ws.Range("I97").Value = 2006.2222
$ws.Range("K97").Value = 2006.2222
$ws.Range("M97").Value = -1510.2222
$ws.Range("H132").Value = 11497340
$ws.Range("I132").Value = 12823466
$ws.Range("K132").Value = 38470398
$ws.Range("M132").Value = -38467868
$ws.Range("H140").Value = 80666.336
$ws.Range("J140").Value = 80666.336
$ws.Range("L140").Value = 80666.336
$ws.Range("N140").Value = -91026.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2884.6316
$ws.Range("I22").Value = 1421.75
$ws.Range("J22").Value = 3274.7334
$ws.Range("K22").Value = 1421.75
$ws.Range("L22").Value = 3274.7334
$ws.Range("M22").Value = -1126.75
$ws.Range("N22").Value = -3864.7334
$ws.Range("H27").Value = 2884.6316
$ws.Range("I27").Value = 1421.75
$ws.Range("J27").Value = 3274.7334
$ws.Range("K27").Value = 1421.75
$ws.Range("L27").Value = 3274.7334
$ws.Range("M27").Value = -1314.75
$ws.Range("N27").Value = -3488.7334
$ws.Range("H46").Value = 8987.706
$ws.Range("I46").Value = 4400
$ws.Range("K46").Value = 4400
$ws.Range("M46").Value = -4212
$ws.Range("H61").Value = 5324.4375
$ws.Range("I61").Value = 5998.5386
$ws.Range("K61").Value = 5998.5386
$ws.Range("M61").Value = -5796.5386
$ws.Range("H93").Value = 1860.5
$ws.Range("I93").Value = 2314.3333
$ws.Range("K93").Value = 2314.3333
$ws.Range("M93").Value = -1066.3333
$ws.Range("H98").Value = 71766.336
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 71766.336
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 71766.336
$ws.Range("M98").Value = ""
$ws.Range("N98").Value = -77756.336
$ws.Range("H113").Value = 5324.4375
$ws.Range("I113").Value = 5998.5386
$ws.Range("K113").Value = 5998.5386
$ws.Range("M113").Value = -3828.5386
$ws.Range("H136").Value = 4081.8
$ws.Range("I136").Value = 3850.7827
$ws.Range("J136").Value = 4840.857
$ws.Range("K136").Value = 11552.3481
$ws.Range("L136").Value = 14522.571
$ws.Range("M136").Value = -9002.348100000001
$ws.Range("N136").Value = -19622.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 52857
$ws.Range("J97").Value = 52857
$ws.Range("L97").Value = 52857
$ws.Range("N97").Value = -54839
$ws.Range("H107").Value = 611.7586
$ws.Range("I107").Value = 317.11765
$ws.Range("J107").Value = 1029.1666
$ws.Range("K107").Value = 951.3529500000001
$ws.Range("L107").Value = 3087.4998
$ws.Range("M107").Value = 968.6470499999999
$ws.Range("N107").Value = -6927.4998
$ws.Range("H128").Value = 80000
$ws.Range("J128").Value = 80000
$ws.Range("L128").Value = 80000
$ws.Range("N128").Value = -89960
$ws.Range("H136").Value = 12604.611
$ws.Range("I136").Value = 19964
$ws.Range("K136").Value = 59892
$ws.Range("M136").Value = -57342
